# Applies the "Add files via upload" revision to Tabela.xlsx:
#  - a few existing table cells get new/changed text
#  - row 4 gains a wrapped "E4" value (and grows to a two-line row)
#  - row 12's empty E12 cell becomes a wrap-text styled cell
#  - two brand-new rows (16 and 17) are appended to the User-Story table
#  - the active selection moves to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing-cell content tweaks -----------------------------------------

# Row 4: Director's "so that" cell was blank, now documents the goal.
$ws.Range("E4").Value = "Upewnienie się, że każdy członek zespołu wykonuje  swoje zadanie"

# Row 5: role changes from the generic "Employee" to "Safety specialist".
$ws.Range("C5").Value = "Safety specialist"

# Row 10: the placeholder "Netcode" task is replaced with the real task.
$ws.Range("D10").Value = "Budowa sieci w budynku"

# Row 15: role changes from "Employee" to "Server Manager".
$ws.Range("C15").Value = "Server Manager"

# --- formatting tweaks on existing cells -----------------------------------

# Header cell E2 gains wrap text (matches the rest of the header row).
$ws.Range("E2").WrapText = $true

# Row 4 becomes a two-line row once E4 carries real wrapped text.
$ws.Rows("4").RowHeight = 30

# Row 12's empty E12 cell switches from a plain bordered cell to a
# wrap-text bordered cell, matching its D12 neighbour.
$ws.Range("E12").WrapText = $true

# --- new rows ---------------------------------------------------------------

# Row 16: new User Story #14.
$ws.Range("B16").Value = 14
$ws.Range("B16").Borders.LineStyle = 1

$ws.Range("C16").Value = "Employee"
$ws.Range("C16").Borders.LineStyle = 1

$ws.Range("D16").Value = "Praca nad API"
$ws.Range("D16").Borders.LineStyle = 1
$ws.Range("D16").WrapText = $true

$ws.Range("E16").Value = "Opracowanie API"
$ws.Range("E16").Borders.LineStyle = 1
$ws.Range("E16").WrapText = $true

# Row 17: new User Story #15.
$ws.Range("B17").Value = 15
$ws.Range("B17").Borders.LineStyle = 1

$ws.Range("C17").Value = "Tester"
$ws.Range("C17").Borders.LineStyle = 1

$ws.Range("D17").Value = "Testowanie aplikacji"
$ws.Range("D17").Borders.LineStyle = 1
$ws.Range("D17").WrapText = $true

$ws.Range("E17").Value = "Upewnienie się, że aplikacja przechodzi wszystkie wymagane testy"
$ws.Range("E17").Borders.LineStyle = 1
$ws.Range("E17").WrapText = $true

$ws.Rows("17").RowHeight = 30

# --- selection --------------------------------------------------------------

$ws.Range("C9").Select() | Out-Null
